$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1223
$ws.Range("F5").Value = 1431
$ws.Range("F6").Value = 1728
$ws.Range("F7").Value = 6259
$ws.Range("F8").Value = 129
$ws.Range("F9").Value = 1851
$ws.Range("F15").Value = 38
$ws.Range("F16").Value = 7063
$ws.Range("F17").Value = 131
$ws.Range("F21").Value = 1721
$ws.Range("F25").Value = 167
$ws.Range("F26").Value = 1634
$ws.Range("F27").Value = 776
$ws.Range("F28").Value = 330
$ws.Range("F33").Value = 3902

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 203
$ws.Range("F23").Value = 3
$ws.Range("F25").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 665
$ws.Range("F5").Value = 252

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 665
$ws.Range("F5").Value = 1223
$ws.Range("F10").Value = 1431
$ws.Range("F11").Value = 252
$ws.Range("F12").Value = 1728
$ws.Range("F13").Value = 6259
$ws.Range("F14").Value = 129
$ws.Range("F15").Value = 1851
$ws.Range("F24").Value = 38
$ws.Range("F25").Value = 7063
$ws.Range("F26").Value = 131
$ws.Range("F30").Value = 1721
$ws.Range("F34").Value = 1634
$ws.Range("F36").Value = 330
$ws.Range("F43").Value = 3902
$ws.Range("F47").Value = 3
$ws.Range("F49").Value = 3

